$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,4,5,6,7,8,9,10) were reshuffled (row 3 unchanged).
# Apply the target values for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg).

$rows = @{
    2  = @{ D = 44750; J = 140; K = 19000; L = 20000; M = 19571; P = 1305 }
    4  = @{ D = 44749; J = 90;  K = 17000; L = 18000; M = 17556; P = 1170 }
    5  = @{ D = 45084; J = 90;  K = 22000; L = 23000; M = 22556; P = 1504 }
    6  = @{ D = 44839; J = 50;  K = 15000; L = 16000; M = 15600; P = 1040 }
    7  = @{ D = 45141; J = 50;  K = 8500;  L = 9000;  M = 8800;  P = 587  }
    8  = @{ D = 45091; J = 40;  K = 20000; L = 22000; M = 21000; P = 1400 }
    9  = @{ D = 45063; J = 40;  K = 21000; L = 22000; M = 21500; P = 1433 }
    10 = @{ D = 45133; J = 50;  K = 22000; L = 22000; M = 22000; P = 1467 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
